# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# columns for the 5b29b475-... file rows on the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: rows 2 and 3 (5b29b475-... handoff/handback timestamps) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-25 08:42:27"
$wsZh.Range("H2").Value = "2016-03-25 08:42:54"
$wsZh.Range("E3").Value = "2016-03-25 08:42:27"
$wsZh.Range("H3").Value = "2016-03-25 08:42:54"

# --- de-de sheet: rows 3 and 4 (5b29b475-... handoff/handback timestamps) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-25 08:42:32"
$wsDe.Range("H3").Value = "2016-03-25 08:43:02"
$wsDe.Range("E4").Value = "2016-03-25 08:42:32"
$wsDe.Range("H4").Value = "2016-03-25 08:43:02"
